$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields ---
$ws.Range("C2").Value = "Hartmut"

# Card number B3 must stay textual (it's a 16-digit number that should not
# be coerced to a numeric / scientific-notation value).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 24.03.2025"

# --- Row 6 ---
$ws.Range("B6").Value = "25.03."
$ws.Range("C6").Value = "26.03."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 18677314"
$ws.Range("E6").Value = "38,03-"

# --- Row 7 ---
$ws.Range("B7").Value = "29.03."
$ws.Range("C7").Value = "30.03."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-76895362"
$ws.Range("E7").Value = "57,51-"

# --- Row 8 ---
$ws.Range("B8").Value = "02.04."
$ws.Range("C8").Value = "03.04."
$ws.Range("D8").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E8").Value = "46,60-"

# --- Row 9 ---
$ws.Range("B9").Value = "04.04."
$ws.Range("C9").Value = "05.04."
$ws.Range("D9").Value = "KARTENZ./04.04 LIDL RO"
$ws.Range("E9").Value = "116,62-"

# --- Row 10: transaction removed, cells blanked ---
$ws.Range("B10:D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

# --- Row 11: transaction removed, cells blanked ---
$ws.Range("B11:D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 07.04.2025"
$ws.Range("E12").Value = "258,76-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 12.04.2025"
